$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '21.741.32'
$ws.Range('E2').Value = '  +5.99%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.573.55'

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9983'

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.9762'
$ws.Range('E5').Value = '  +1.97%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '285.09'
$ws.Range('E6').Value = '  +2.99%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3677'
$ws.Range('E7').Value = '  +0.63%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3270'
$ws.Range('E8').Value = '  +6.86%  '

$ws.Range('E9').Value = '  +7.49%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '41.15'
$ws.Range('E10').Value = '  +3.70%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07050'
$ws.Range('E11').Value = '  +6.68%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.9940'
$ws.Range('E12').Value = '  -0.81%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '20.19'
$ws.Range('E13').Value = '  +11.50%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.813'
$ws.Range('E14').Value = '  +6.54%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.477'
$ws.Range('E15').Value = '  +4.85%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.00001070'
$ws.Range('E16').Value = '  +4.08%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.9749'
$ws.Range('E17').Value = '  +1.22%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '1.562.42'
$ws.Range('E18').Value = '  +5.93%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06177'
$ws.Range('E19').Value = '  +4.67%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '73.91'
$ws.Range('E20').Value = '  +7.00%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '15.96'
$ws.Range('E21').Value = '  +10.25%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.835'
$ws.Range('E22').Value = '  +6.78%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '11.58'
$ws.Range('E23').Value = '  +5.03%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '21.727.13'
$ws.Range('E24').Value = '  +5.64%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.320'
$ws.Range('E25').Value = '  +2.65%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.412'
$ws.Range('E26').Value = '  +13.86%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '148.37'
$ws.Range('E27').Value = '  +5.23%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.12'
$ws.Range('E28').Value = '  +5.52%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.736.36'
$ws.Range('E29').Value = '  +6.41%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '120.33'
$ws.Range('E30').Value = '  +5.98%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.060'
$ws.Range('E31').Value = '  +2.81%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.9021'

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.387'
$ws.Range('E33').Value = '  +8.64%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.08166'
$ws.Range('E34').Value = '  +2.81%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.581'
$ws.Range('E35').Value = '  +3.19%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.112'
$ws.Range('E36').Value = '  +8.12%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '11.53'
$ws.Range('E37').Value = '  +10.63%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.230'
$ws.Range('E38').Value = '  +0.41%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.06015'
$ws.Range('E39').Value = '  +4.31%  '

$ws.Range('B40').Value = 'VeChain'
$ws.Range('C40').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.02159'
$ws.Range('E40').Value = '  +6.10%  '

$ws.Range('B41').Value = 'FraxShare'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '8.106'
$ws.Range('E41').Value = '  +6.83%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1999'
$ws.Range('E42').Value = '  +6.48%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.9742'
$ws.Range('E43').Value = '  +1.89%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.5709'
$ws.Range('E44').Value = '  +8.24%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '12.74'
$ws.Range('E45').Value = '  +6.18%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.605'
$ws.Range('E46').Value = '  +2.92%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5622'
$ws.Range('E47').Value = '  +8.37%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '124.44'
$ws.Range('E48').Value = '  +5.99%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.921'
$ws.Range('E49').Value = '  +8.12%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06727'
$ws.Range('E50').Value = '  +4.21%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '71.45'
$ws.Range('E51').Value = '  +6.46%  '
